$wb = $excel.ActiveWorkbook

# --- Rename the second worksheet ("Bioreactor cultivation conditio" -> "bioreactor__cultivation") ---
$ws = $wb.Worksheets.Item(2)
$ws.Name = "bioreactor__cultivation"

# --- Update the "Table" name cell on the isa_template sheet to match the new sheet/table name ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B7").Value = "bioreactor__cultivation"

# --- Fill in data row values on the renamed sheet ---
# Input [Sample Name] for the data row
$ws.Range("A2").Value = "S1"

# Parameter [Oxygen] value changes from "20,95" to "20" (kept as text, matching
# the existing "20" used elsewhere in the row, e.g. Parameter [Carbon Dioxide])
$ws.Range("AI2").NumberFormat = "@"
$ws.Range("AI2").Value = "20"
$ws.Range("AI2").ClearFormats()

# Output [Data] for the data row
$ws.Range("BF2").Value = "O1"
